# Make fixture more complex:
#  - add a value to C1
#  - add two new rows of data (rows 3 and 4)
#  - merge B3:B4
#  - update the active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cell on row 1
$ws.Range("C1").Value = 3

# New row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 4
$ws.Range("C3").Value = 3

# New row 4
$ws.Range("A4").Value = 1
$ws.Range("C4").Value = 3

# Merge B3:B4 (value 4 lives in the merged cell)
$ws.Range("B3:B4").Merge()

# Move the selection, matching the recorded sheet view state
$ws.Range("J23").Select()
